$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$g = $s.Shapes.Item(1)

$sh = $g.GroupItems.Item("rc4")
$sh.Left = 90.205512
$sh.Top = 452.181260
$sh.Width = 65.721732
$sh.Height = 27.253858

$sh = $g.GroupItems.Item("rc5")
$sh.Left = 163.229685
$sh.Top = 403.124173
$sh.Width = 65.721732
$sh.Height = 76.310948

$sh = $g.GroupItems.Item("rc6")
$sh.Left = 236.253945
$sh.Top = 375.870236
$sh.Width = 65.721732
$sh.Height = 103.564885

$sh = $g.GroupItems.Item("rc7")
$sh.Left = 309.278122
$sh.Top = 348.616299
$sh.Width = 65.721732
$sh.Height = 130.818825

$sh = $g.GroupItems.Item("rc8")
$sh.Left = 382.302292
$sh.Top = 239.600551
$sh.Width = 65.721732
$sh.Height = 239.834572

$sh = $g.GroupItems.Item("rc10")
$sh.Left = 528.350739
$sh.Top = 468.533622
$sh.Width = 65.721732
$sh.Height = 10.901497

$sh = $g.GroupItems.Item("tx11")
$sh.Left = 119.110394
$sh.Top = 406.142205
$sh.Width = 7.911969
$sh.Height = 10.218189
$sh.TextFrame.TextRange.Text = "5"

$sh = $g.GroupItems.Item("tx12")
$sh.Left = 108.048110
$sh.Top = 423.495354
$sh.TextFrame.TextRange.Text = "(3%)"

$sh = $g.GroupItems.Item("tx13")
$sh.Left = 188.178583
$sh.Top = 357.078202

$sh = $g.GroupItems.Item("tx14")
$sh.Left = 181.072289
$sh.Top = 374.438279

$sh = $g.GroupItems.Item("tx15")
$sh.Left = 261.202756
$sh.Top = 329.650630

$sh = $g.GroupItems.Item("tx16")
$sh.Left = 250.140481
$sh.Top = 347.184342

$sh = $g.GroupItems.Item("tx17")
$sh.Left = 334.226944
$sh.Top = 302.570328
$sh.Width = 15.824017
$sh.Height = 10.225198
$sh.TextFrame.TextRange.Text = "24"

$sh = $g.GroupItems.Item("tx18")
$sh.Left = 323.164724
$sh.Top = 319.930405
$sh.TextFrame.TextRange.Text = "(15%)"

$sh = $g.GroupItems.Item("tx19")
$sh.Left = 407.251181
$sh.Top = 193.596220
$sh.Width = 15.824017
$sh.Height = 10.183465
$sh.TextFrame.TextRange.Text = "44"

$sh = $g.GroupItems.Item("tx20")
$sh.Left = 396.188898
$sh.Top = 210.914650
$sh.TextFrame.TextRange.Text = "(28%)"

$sh = $g.GroupItems.Item("tx21")
$sh.Left = 480.275354
$sh.Top = 155.225433
$sh.Width = 15.824017
$sh.Height = 10.398819
$sh.TextFrame.TextRange.Text = "51"

$sh = $g.GroupItems.Item("tx22")
$sh.TextFrame.TextRange.Text = "(32%)"

$sh = $g.GroupItems.Item("tx23")
$sh.Left = 557.255591
$sh.Top = 422.487638
$sh.Width = 7.911969
$sh.Height = 10.225198
$sh.TextFrame.TextRange.Text = "2"

$sh = $g.GroupItems.Item("tx24")
$sh.Left = 546.193329
$sh.Top = 439.847717
$sh.TextFrame.TextRange.Text = "(1%)"

$sh = $g.GroupItems.Item("tx27")
$sh.Left = 59.428504
$sh.Top = 365.942835

$sh = $g.GroupItems.Item("tx28")
$sh.Left = 59.428504
$sh.Top = 256.927087

$sh = $g.GroupItems.Item("tx29")
$sh.Left = 59.428504
$sh.Top = 147.911339

$sh = $g.GroupItems.Item("pl31")
$sh.Left = 75.764961
$sh.Top = 370.419449

$sh = $g.GroupItems.Item("pl32")
$sh.Left = 75.764961
$sh.Top = 261.403702

$sh = $g.GroupItems.Item("pl33")
$sh.Left = 75.764961
$sh.Top = 152.387953

$sh = $g.GroupItems.Item("tx51")
$sh.TextFrame.TextRange.Text = "diagnostic de trouble psychotique? (n = 159)"
